# Update the "K" column (column G) values in the save-data sheet.
# These are pre-computed strikeout-count values (K) that replace the
# previous Strike# values, as part of regenerating save_data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 9
    4  = 8
    5  = 7
    6  = 3
    7  = 4
    8  = 4
    9  = 2
    10 = 7
    11 = 5
    12 = 4
    13 = 7
    14 = 5
    15 = 6
    16 = 5
    17 = 6
    18 = 5
    19 = 5
    20 = 7
    21 = 6
    22 = 5
    23 = 5
    24 = 6
    25 = 11
    26 = 2
    27 = 6
    28 = 6
    29 = 6
    30 = 4
    31 = 5
    32 = 6
    33 = 3
    34 = 11
    35 = 4
    36 = 2
    37 = 2
    38 = 3
    39 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
